$wb = $excel.ActiveWorkbook

# Rename sheets (task order ids updated)
$wb.Worksheets.Item(1).Name = "GNG_TO-16512556209355602"
$wb.Worksheets.Item(2).Name = "NB_TO-16512556253085067"
$wb.Worksheets.Item(3).Name = "RS_TO-16512556253103416"
$wb.Worksheets.Item(4).Name = "TOL_TO-1651255625356043"
$wb.Worksheets.Item(5).Name = "vSAT_TO-1651255625433689"

# Sheet1 (GNG)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16512556208945637.csv"
$ws1.Range("B3").Value = "GNG_stims-16512556209175577.csv"
$ws1.Range("B4").Value = "go_stims-16512556209195583.csv"
$ws1.Range("B5").Value = "GNG_stims-16512556209335601.csv"

# Sheet2 (NB)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-16512556219568968.csv"
$ws2.Range("B3").Value = "TB-16512556251810222.csv"
$ws2.Range("B4").Value = "TB-16512556252882862.csv"
$ws2.Range("B5").Value = "OB-1651255623047286.csv"
$ws2.Range("B6").Value = "ZB-match_8-1651255621563529.csv"
$ws2.Range("B7").Value = "ZB-match_7-1651255621842907.csv"
$ws2.Range("B8").Value = "TB-165125562332639.csv"
$ws2.Range("B9").Value = "ZB-match_9-16512556212512975.csv"
$ws2.Range("B10").Value = "OB-16512556228870242.csv"

# Sheet3 (RS)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes open"
$ws3.Range("B3").Value = "eyes closed"

# Sheet4 (TOL)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16512556253243928.csv"
$ws4.Range("B3").Value = "ZM_stims-16512556253133402.csv"
$ws4.Range("B4").Value = "MM_stims-16512556253398125.csv"
$ws4.Range("B5").Value = "ZM_stims-16512556253253915.csv"
$ws4.Range("B6").Value = "MM_stims-1651255625356043.csv"
$ws4.Range("B7").Value = "ZM_stims-16512556253408103.csv"

# Sheet5 (vSAT)
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-16512556253616533.csv"
$ws5.Range("B3").Value = "vSAT_stims-16512556254036276.csv"
$ws5.Range("B4").Value = "vSAT_stims-1651255625418362.csv"
$ws5.Range("B5").Value = "SAT_stims-16512556253884137.csv"
